# Update workbook with new data rows (aggiornamento fino a 13/03)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data to append below the existing last row (row 251)
$newRows = @(
    @{ Row = 252; A = 44326; B = 0; C = 2; D = 81.59934720522236 },
    @{ Row = 253; A = 44327; B = 0; C = 1; D = 40.79967360261118 },
    @{ Row = 254; A = 44328; B = 0; C = 1; D = 40.79967360261118 },
    @{ Row = 255; A = 44329; B = 0; C = 1; D = 40.79967360261118 }
)

foreach ($r in $newRows) {
    $rowIndex = $r.Row
    $prevRow = $rowIndex - 1

    # Copy the formatting of the date cell in the previous (last existing) row
    # so the new date cell keeps the same style (s="2": date format + border).
    $ws.Range("A$prevRow").Copy()
    $ws.Range("A$rowIndex").PasteSpecial(-4122)

    $ws.Cells.Item($rowIndex, 1).Value = $r.A
    $ws.Cells.Item($rowIndex, 2).Value = $r.B
    $ws.Cells.Item($rowIndex, 3).Value = $r.C
    $ws.Cells.Item($rowIndex, 4).Value = $r.D
}

$excel.CutCopyMode = $false
